# Scheduled-runner refresh: write back updated market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) for the affected Leve rows
# across each crafting-class sheet's Table_<CLASS> listing.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3335.8235
$ws.Range("I19").Value = 6231.4116
$ws.Range("K19").Value = 6231.4116
$ws.Range("M19").Value = -6056.4116
$ws.Range("H28").Value = 581.48
$ws.Range("I28").Value = 598.86365
$ws.Range("J28").Value = 454
$ws.Range("K28").Value = 598.86365
$ws.Range("L28").Value = 454
$ws.Range("M28").Value = -113.86365
$ws.Range("N28").Value = -1424
$ws.Range("H64").Value = 5860.857
$ws.Range("I64").Value = 4879.6
$ws.Range("J64").Value = 6406
$ws.Range("K64").Value = 4879.6
$ws.Range("L64").Value = 6406
$ws.Range("M64").Value = -4631.6
$ws.Range("N64").Value = -6902
$ws.Range("H67").Value = 5860.857
$ws.Range("I67").Value = 4879.6
$ws.Range("J67").Value = 6406
$ws.Range("K67").Value = 4879.6
$ws.Range("L67").Value = 6406
$ws.Range("M67").Value = -4021.6
$ws.Range("N67").Value = -8122
$ws.Range("H74").Value = 5593.357
$ws.Range("I74").Value = 5588
$ws.Range("J74").Value = 5596.3335
$ws.Range("K74").Value = 5588
$ws.Range("L74").Value = 5596.3335
$ws.Range("M74").Value = -4652
$ws.Range("N74").Value = -7468.3335
$ws.Range("H76").Value = 6160.203
$ws.Range("I76").Value = 6532.396
$ws.Range("K76").Value = 6532.396
$ws.Range("M76").Value = -6217.396
$ws.Range("H77").Value = 5593.357
$ws.Range("I77").Value = 5588
$ws.Range("J77").Value = 5596.3335
$ws.Range("K77").Value = 27940
$ws.Range("L77").Value = 27981.6675
$ws.Range("M77").Value = -23260
$ws.Range("N77").Value = -37341.6675
$ws.Range("H79").Value = 6160.203
$ws.Range("I79").Value = 6532.396
$ws.Range("K79").Value = 6532.396
$ws.Range("M79").Value = -5440.396

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 493.34784
$ws.Range("I2").Value = 444.6
$ws.Range("K2").Value = 444.6
$ws.Range("M2").Value = -331.6
$ws.Range("H63").Value = 4901.56
$ws.Range("I63").Value = 5361.65
$ws.Range("J63").Value = 3061.2
$ws.Range("K63").Value = 5361.65
$ws.Range("L63").Value = 3061.2
$ws.Range("M63").Value = -4675.65
$ws.Range("N63").Value = -4433.2
$ws.Range("H66").Value = 4901.56
$ws.Range("I66").Value = 5361.65
$ws.Range("J66").Value = 3061.2
$ws.Range("K66").Value = 26808.25
$ws.Range("L66").Value = 15306
$ws.Range("M66").Value = -23376.25
$ws.Range("N66").Value = -22170
$ws.Range("H116").Value = 493.34784
$ws.Range("I116").Value = 444.6
$ws.Range("K116").Value = 444.6
$ws.Range("M116").Value = 1849.4
$ws.Range("H122").Value = 2185.1936
$ws.Range("I122").Value = 2296.1738
$ws.Range("J122").Value = 1866.125
$ws.Range("K122").Value = 6888.5214
$ws.Range("L122").Value = 5598.375
$ws.Range("M122").Value = -4438.5214
$ws.Range("N122").Value = -10498.375
$ws.Range("H135").Value = 22737.182
$ws.Range("J135").Value = 22737.182
$ws.Range("L135").Value = 22737.182
$ws.Range("N135").Value = -32877.182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 493.34784
$ws.Range("I3").Value = 444.6
$ws.Range("K3").Value = 444.6
$ws.Range("M3").Value = -330.6
$ws.Range("H105").Value = 2588.6978
$ws.Range("I105").Value = 2416.75
$ws.Range("K105").Value = 2416.75
$ws.Range("M105").Value = -669.75
$ws.Range("H134").Value = 57658.5
$ws.Range("I134").Value = 88547.5
$ws.Range("J134").Value = 11325
$ws.Range("K134").Value = 265642.5
$ws.Range("L134").Value = 33975
$ws.Range("M134").Value = -263107.5
$ws.Range("N134").Value = -39045

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1269.5883
$ws.Range("I16").Value = 1508.75
$ws.Range("J16").Value = 695.6
$ws.Range("K16").Value = 1508.75
$ws.Range("L16").Value = 695.6
$ws.Range("M16").Value = -1221.75
$ws.Range("N16").Value = -1269.6
$ws.Range("H31").Value = 3054.418
$ws.Range("I31").Value = 2114.0605
$ws.Range("K31").Value = 2114.0605
$ws.Range("M31").Value = -1819.0605
$ws.Range("H34").Value = 3054.418
$ws.Range("I34").Value = 2114.0605
$ws.Range("K34").Value = 2114.0605
$ws.Range("M34").Value = -1912.0605
$ws.Range("H62").Value = 3221.976
$ws.Range("I62").Value = 2980.2104
$ws.Range("J62").Value = 3421.6956
$ws.Range("K62").Value = 2980.2104
$ws.Range("L62").Value = 3421.6956
$ws.Range("M62").Value = -2356.2104
$ws.Range("N62").Value = -4669.6956
$ws.Range("H65").Value = 3221.976
$ws.Range("I65").Value = 2980.2104
$ws.Range("J65").Value = 3421.6956
$ws.Range("K65").Value = 14901.052
$ws.Range("L65").Value = 17108.478
$ws.Range("M65").Value = -11781.052
$ws.Range("N65").Value = -23348.478
$ws.Range("H99").Value = 75249.64
$ws.Range("I99").Value = 114201.664
$ws.Range("J99").Value = 5136
$ws.Range("K99").Value = 114201.664
$ws.Range("L99").Value = 5136
$ws.Range("M99").Value = -112703.664
$ws.Range("N99").Value = -8132
$ws.Range("H113").Value = 1269.5883
$ws.Range("I113").Value = 1508.75
$ws.Range("J113").Value = 695.6
$ws.Range("K113").Value = 1508.75
$ws.Range("L113").Value = 695.6
$ws.Range("M113").Value = 661.25
$ws.Range("N113").Value = -5035.6
$ws.Range("H126").Value = 75249.64
$ws.Range("I126").Value = 114201.664
$ws.Range("J126").Value = 5136
$ws.Range("K126").Value = 342604.992
$ws.Range("L126").Value = 15408
$ws.Range("M126").Value = -340134.992
$ws.Range("N126").Value = -20348
$ws.Range("H134").Value = 2127.4644
$ws.Range("I134").Value = 1291
$ws.Range("J134").Value = 2754.8125
$ws.Range("K134").Value = 3873
$ws.Range("L134").Value = 8264.4375
$ws.Range("M134").Value = -1338
$ws.Range("N134").Value = -13334.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 84.04761999999999
$ws.Range("I12").Value = 28.2
$ws.Range("K12").Value = 84.59999999999999
$ws.Range("M12").Value = 88.40000000000001
$ws.Range("H70").Value = 5350
$ws.Range("I70").Value = 900
$ws.Range("J70").Value = 5985.7144
$ws.Range("K70").Value = 2700
$ws.Range("L70").Value = 17957.1432
$ws.Range("M70").Value = -2385
$ws.Range("N70").Value = -18587.1432
$ws.Range("H73").Value = 5350
$ws.Range("I73").Value = 900
$ws.Range("J73").Value = 5985.7144
$ws.Range("K73").Value = 2700
$ws.Range("L73").Value = 17957.1432
$ws.Range("M73").Value = -1608
$ws.Range("N73").Value = -20141.1432
$ws.Range("H76").Value = 5484.5
$ws.Range("I76").Value = 998
$ws.Range("K76").Value = 2994
$ws.Range("M76").Value = -2611
$ws.Range("H79").Value = 5484.5
$ws.Range("I79").Value = 998
$ws.Range("K79").Value = 2994
$ws.Range("M79").Value = -1668
$ws.Range("H113").Value = 501.54166
$ws.Range("I113").Value = 481
$ws.Range("J113").Value = 520.9729599999999
$ws.Range("K113").Value = 1443
$ws.Range("L113").Value = 1562.91888
$ws.Range("M113").Value = 727
$ws.Range("N113").Value = -5902.918879999999
$ws.Range("H131").Value = 862.9838999999999
$ws.Range("I131").Value = 469.7143
$ws.Range("J131").Value = 913.0364
$ws.Range("K131").Value = 1409.1429
$ws.Range("L131").Value = 2739.1092
$ws.Range("M131").Value = 3630.8571
$ws.Range("N131").Value = -12819.1092
$ws.Range("H132").Value = 2976.6
$ws.Range("I132").Value = 1039.5
$ws.Range("J132").Value = 3622.3
$ws.Range("K132").Value = 9355.5
$ws.Range("L132").Value = 32600.7
$ws.Range("M132").Value = -6825.5
$ws.Range("N132").Value = -37660.7
$ws.Range("H134").Value = 2717.4517
$ws.Range("I134").Value = 1790.4286
$ws.Range("J134").Value = 3480.8823
$ws.Range("K134").Value = 5371.2858
$ws.Range("L134").Value = 10442.6469
$ws.Range("M134").Value = -301.2857999999997
$ws.Range("N134").Value = -20582.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9038.348
$ws.Range("I70").Value = 3887.6875
$ws.Range("K70").Value = 3887.6875
$ws.Range("M70").Value = -3617.6875
$ws.Range("H73").Value = 9038.348
$ws.Range("I73").Value = 3887.6875
$ws.Range("K73").Value = 3887.6875
$ws.Range("M73").Value = -2951.6875
$ws.Range("H80").Value = 4230.227
$ws.Range("I80").Value = 4639.1177
$ws.Range("K80").Value = 4639.1177
$ws.Range("M80").Value = -3641.1177
$ws.Range("H83").Value = 4230.227
$ws.Range("I83").Value = 4639.1177
$ws.Range("K83").Value = 23195.5885
$ws.Range("M83").Value = -18203.5885
$ws.Range("H107").Value = 4699.4346
$ws.Range("I107").Value = 6423.75
$ws.Range("J107").Value = 758.1429000000001
$ws.Range("K107").Value = 6423.75
$ws.Range("L107").Value = 758.1429000000001
$ws.Range("M107").Value = -4503.75
$ws.Range("N107").Value = -4598.1429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 311.90475
$ws.Range("I107").Value = 324.75
$ws.Range("J107").Value = 294.77777
$ws.Range("K107").Value = 974.25
$ws.Range("L107").Value = 884.33331
$ws.Range("M107").Value = 945.75
$ws.Range("N107").Value = -4724.33331
$ws.Range("H122").Value = 29593.742
$ws.Range("I122").Value = 34367.965
$ws.Range("J122").Value = 948.4
$ws.Range("K122").Value = 103103.895
$ws.Range("L122").Value = 2845.2
$ws.Range("M122").Value = -100653.895
$ws.Range("N122").Value = -7745.2
$ws.Range("H136").Value = 32260866
$ws.Range("I136").Value = 71431060
$ws.Range("J136").Value = 3064.7058
$ws.Range("K136").Value = 214293180
$ws.Range("L136").Value = 9194.117400000001
$ws.Range("M136").Value = -214290630
$ws.Range("N136").Value = -14294.1174
